# The commit swaps the two theme parts in the deck: the Slide Master
# (ppt/theme/theme1.xml) switches from the custom "Integral" / "Red Violet"
# colour scheme to the stock "Office Theme" / "Office" colour scheme (the
# scheme previously used only by the Notes Master). Font scheme and format
# scheme (fills/lines/effects) are identical between the two themes, so the
# only real content change is the 12-slot theme colour scheme.
#
# We reach the Slide Master's theme colours through the PowerPoint object
# model's ThemeColorScheme (exposed on a Slide, but it edits the single
# shared master theme used by every slide in this deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> theme slot (verified empirically against the underlying OOXML):
#  1 dk1   2 lt1   3 dk2   4 lt2
#  5 accent1  6 accent2  7 accent3  8 accent4
#  9 accent5  10 accent6  11 hlink  12 folHlink
#
# Target values = the stock Office theme colour scheme.
# PowerPoint's RGB property is a VBA-style packed integer: R + G*256 + B*65536.
$tcs.Colors(1).RGB  = 0x000000   # dk1      srgb 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      srgb FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      srgb 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      srgb E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  srgb 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  srgb ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  srgb A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  srgb FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  srgb 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  srgb 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    srgb 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink srgb 954F72
